# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 3.32 = 12721.13 pesos", "1000 Bs = 3.32 = 12757.48 pesos")
$text = $text.Replace("12721.13 pesos = 3.31 = 971.31 Bs", "12757.48 pesos = 3.31 = 981.88 Bs")
$cell.Value = $text

# --- tasas: update numeric rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3840
$wsTasas.Range("N12").Value = 3855
$wsTasas.Range("O12").Value = 296.701
